# Update CV results (B2:E19) for base1 gridsearch xgb results - ejecucion 1
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 161.9419643878937
$ws.Range("C2").Value = 0.3660026642539724
$ws.Range("D2").Value = 1.6385995388031
$ws.Range("E2").Value = 0.03638119533530275
$ws.Range("B3").Value = 330.6253034591675
$ws.Range("C3").Value = 2.442393107128127
$ws.Range("D3").Value = 1.691832399368286
$ws.Range("E3").Value = 0.04152285251884878
$ws.Range("B4").Value = 657.8522876739502
$ws.Range("C4").Value = 3.657536074083557
$ws.Range("D4").Value = 1.821935796737671
$ws.Range("E4").Value = 0.1585939345128875
$ws.Range("B5").Value = 169.9950989246368
$ws.Range("C5").Value = 2.36209444418793
$ws.Range("D5").Value = 1.69543628692627
$ws.Range("E5").Value = 0.06191757124431183
$ws.Range("B6").Value = 331.1381764411926
$ws.Range("C6").Value = 0.9196429678536759
$ws.Range("D6").Value = 1.968860912322998
$ws.Range("E6").Value = 0.2356524003722107
$ws.Range("B7").Value = 660.6188005924225
$ws.Range("C7").Value = 1.924461351001752
$ws.Range("D7").Value = 1.837553977966309
$ws.Range("E7").Value = 0.2177630975062967
$ws.Range("B8").Value = 172.5016034126282
$ws.Range("C8").Value = 1.709480852340438
$ws.Range("D8").Value = 1.78572883605957
$ws.Range("E8").Value = 0.2291547222481176
$ws.Range("B9").Value = 332.5146837234497
$ws.Range("C9").Value = 1.892363378821995
$ws.Range("D9").Value = 1.850902271270752
$ws.Range("E9").Value = 0.1803837276169671
$ws.Range("B10").Value = 659.1947891235352
$ws.Range("C10").Value = 4.391378112339347
$ws.Range("D10").Value = 1.745503616333008
$ws.Range("E10").Value = 0.1495570587827174
$ws.Range("B11").Value = 273.8755764961243
$ws.Range("C11").Value = 2.014467950726404
$ws.Range("D11").Value = 1.680195999145508
$ws.Range("E11").Value = 0.0613857398344251
$ws.Range("B12").Value = 533.17073097229
$ws.Range("C12").Value = 1.15460117136436
$ws.Range("D12").Value = 1.731850147247314
$ws.Range("E12").Value = 0.1568064772251539
$ws.Range("B13").Value = 1063.094308376312
$ws.Range("C13").Value = 3.597990834486031
$ws.Range("D13").Value = 1.701994132995605
$ws.Range("E13").Value = 0.1143396499592568
$ws.Range("B14").Value = 271.365793466568
$ws.Range("C14").Value = 1.47232186345623
$ws.Range("D14").Value = 1.704097032546997
$ws.Range("E14").Value = 0.09460969435252703
$ws.Range("B15").Value = 532.5155857086181
$ws.Range("C15").Value = 0.9439214851462323
$ws.Range("D15").Value = 1.696440029144287
$ws.Range("E15").Value = 0.1631662953975631
$ws.Range("B16").Value = 1056.032032823563
$ws.Range("C16").Value = 1.300196273474587
$ws.Range("D16").Value = 1.827653026580811
$ws.Range("E16").Value = 0.1964627566479237
$ws.Range("B17").Value = 269.8764802932739
$ws.Range("C17").Value = 0.848837269330167
$ws.Range("D17").Value = 1.711333084106445
$ws.Range("E17").Value = 0.09356180968766163
$ws.Range("B18").Value = 527.7809926509857
$ws.Range("C18").Value = 5.331695222945193
$ws.Range("D18").Value = 1.821272468566895
$ws.Range("E18").Value = 0.3407503502012928
$ws.Range("B19").Value = 856.1849995136261
$ws.Range("C19").Value = 28.24730491375309
$ws.Range("D19").Value = 1.193986320495605
$ws.Range("E19").Value = 0.262839792029862
